$d = $word.ActiveDocument

# --- Change 1 -------------------------------------------------------------
# "2. %" + "CPU :" (gramStart/gramEnd wrapped) + " percentage of CPU..."
# becomes a single run: "2. %CPU: percentage of CPU..."
$d.Content.Find.Execute(
    "2. %CPU : percentage of CPU resources the process is currently using this is calculated according to the total resources available on the CPU.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2. %CPU: percentage of CPU resources the process is currently using this is calculated according to the total resources available on the CPU.",
    2) | Out-Null

# --- Change 2 -------------------------------------------------------------
# "Memory: The amount of memory the process is currently using shows in
#  different units." becomes the reworded sentence, split across five runs:
#   "Memory: The amount of memory the process is currently using "
#   "is "
#   "show"
#   "n"
#   " in different units."
$d.Content.Find.Execute(
    "Memory: The amount of memory the process is currently using shows in different units.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Memory: The amount of memory the process is currently using is shown in different units.",
    2) | Out-Null

$memRange = $d.Content
$found = $memRange.Find.Execute(
    "Memory: The amount of memory the process is currently using is shown in different units.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $memStart = $memRange.Start
    $memEnd = $memRange.End

    # Offsets (relative to $memStart) at which a new run must begin so the
    # final text lines up with the five pieces listed above.
    $splitOffsets = @(60, 63, 67, 68)

    foreach ($offset in $splitOffsets) {
        $splitPoint = $memStart + $offset
        $tail = $d.Range($splitPoint, $memEnd)
        # Toggling a character property on/off forces a run boundary at
        # $splitPoint without altering the visible formatting.
        $tail.Font.Bold = $true
        $tail.Font.Bold = $false
    }
}
